# Remote tests and Test data changes
# - E/F columns (rows 2-8): "50.07X"/"30.07X" -> "51.05X"/"31.05X"
# - Q column (rows 2-8) and V column formula/result: "test90X@test.com" -> "tests20X@test.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($i = 2; $i -le 8; $i++) {
    $n = $i - 1
    $ws.Range("E$i").Value = "51.05$n"
    $ws.Range("F$i").Value = "31.05$n"
    $ws.Range("Q$i").Value = "tests20$n@test.com"
}

# Update formulas: V2 has its own formula, V3 is the shared-formula master (covers V3:V8)
$ws.Range("V2").Formula = '="tests20" & U2& "@test.com"'
$ws.Range("V3:V8").Formula = '="tests20" & U3& "@test.com"'

# Clear the Q2:Q8 selection left over from the previous edit session
$ws.Range("A1").Select() | Out-Null
